$d = $word.ActiveDocument

# --- Paragraph 1: update indent, add paragraph border, fix ID text, drop trailing space run ---
$p1 = $d.Paragraphs(1)

# Change left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p1.Format.LeftIndent = 225 / 20

# Add a paragraph border (top/left/bottom/right, each with 5pt "space" only, no line)
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Remove the trailing " " run by deleting the paragraph's final space character
# (paragraph text is "**ID__AFFARS_mp_5325_7003_3_topic_3__ID** " + paragraph mark)
$spaceStart = $p1.Range.End - 2
$spaceRange = $d.Range($spaceStart, $spaceStart + 1)
$spaceRange.Delete()

# Fix up the bookmark-style ID text
$d.Content.Find.Execute("**ID__AFFARS_mp_5325_7003_3_topic_3__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_MP5325_7003_3_2__ID**", 2)
